$wb = $excel.ActiveWorkbook

$global = $wb.Worksheets.Item("Global")
$aciclovir = $wb.Worksheets.Item("Aciclovir")

# Create the new "MissingParam" sheet by copying an existing sheet (so it
# inherits the workbook's normal sheet formatting/namespaces) and moving it
# right after "Global" / right before "Aciclovir".
$aciclovir.Copy($null, $global)
$newSheet = $wb.Worksheets.Item("Aciclovir (2)")
$newSheet.Name = "MissingParam"

# Replace the copied data with the definition of the missing parameter.
$newSheet.Range("A1").Value = "Container Path"
$newSheet.Range("B1").Value = "Parameter Name"
$newSheet.Range("C1").Value = "Value"
$newSheet.Range("D1").Value = "Units"

$newSheet.Range("A2").Value = "foo"
$newSheet.Range("B2").Value = "bar"
$newSheet.Range("C2").Value = 2
$newSheet.Range("D2").ClearContents()

$newSheet.Range("C1:C2").NumberFormat = "0.0000"

# Update selection on the Global sheet (no longer the active tab)
$global.Select() | Out-Null
$global.Range("A1:D3").Select() | Out-Null

# Make the new sheet the active / selected tab
$newSheet.Select() | Out-Null
$newSheet.Range("A2:XFD2").Select() | Out-Null
